$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '28.255.37'
$ws.Range("E2").Value = '  +4.31%  '

# Row 3
$ws.Range("D3").Value = '1.786.29'
$ws.Range("E3").Value = '  +0.32%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.002'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.38%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '337.79'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.52%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9978'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.43%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3827'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.07%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3438'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.84%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '47.69'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.63%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.159'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -2.07%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07437'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.08%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '23.49'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +8.88%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.9991'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.43%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.444'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.26%  '

# Row 15
$ws.Range("B15").Value = 'Chainlink'
$ws.Range("C15").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.181'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.22%  '

# Row 16
$ws.Range("B16").Value = 'WrappedEther'
$ws.Range("C16").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D16").Value = '1.786.90'
$ws.Range("E16").Value = '  +0.30%  '

# Row 17
$ws.Range("E17").Value = '  -0.69%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.06663'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.04%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '82.85'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.54%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.9976'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.42%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.53'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.00%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.457'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.62%  '

# Row 23
$ws.Range("D23").Value = '28.276.78'
$ws.Range("E23").Value = '  +4.34%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '12.15'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.45%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.370'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.27%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '20.90'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.99%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.441'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.22%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.425'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -2.76%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '155.00'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.42%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '136.35'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.80%  '

# Row 31
$ws.Range("D31").Value = '1.987.03'
$ws.Range("E31").Value = '  +0.16%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.174'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +3.07%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.968'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.12%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.08874'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +2.44%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '12.85'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.96%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.02441'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +5.03%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.6888'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.22%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '5.354'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.42%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.06378'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.05%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.2187'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.56%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.243'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.60%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.500'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -7.57%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '8.359'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.67%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '14.32'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.33%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.9978'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.40%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.6323'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.10%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.866'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.38%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '133.21'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.76%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.100'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -2.75%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.07471'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +5.30%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.208'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +8.71%  '
